$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new first item ("CORASORE 150MG 20 TAB") above the current
#     first item row (row 7), pushing every item row, the totals row and the
#     footer row down by one row. ---
$ws.Rows("7:7").Insert()

# Copy the formatting + merged-cell layout of the row that used to be row 7
# (now row 8, still the same kind of item row) into the freshly inserted
# (blank) row 7.
$ws.Range("A8:Q8").Copy($ws.Range("A7:Q7"))
$ws.Rows("7:7").RowHeight = 25.5

# Populate the new row 7 with the new item's data. L7/P7 sit on numeric
# number-formats, so a plain text assignment of a numeric-looking string
# ("1", "23.0000") would silently become a real number; flip the format to
# Text for the assignment, then restore the original numeric format code so
# the cell's style stays identical to its neighbours.
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "CORASORE 150MG 20 TAB"
$ws.Range("H7").Value = "6:0"

$fmtL = $ws.Range("L8").NumberFormat
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = "1"
$ws.Range("L7").NumberFormat = $fmtL

$ws.Range("N7").Value = "46.00"

$fmtP = $ws.Range("P8").NumberFormat
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "23.0000"
$ws.Range("P7").NumberFormat = $fmtP

$ws.Range("Q7").Value = "0:1"

# Renumber the items that shifted down by one row (they kept their old
# sequence numbers after the row insert).
$ws.Range("A8").Value = 2
$ws.Range("A9").Value = 3
$ws.Range("A10").Value = 4

# Restore the row heights on the rows that shifted down, matching originals.
$ws.Rows("8:8").RowHeight = 24.75
$ws.Rows("9:9").RowHeight = 25.5
$ws.Rows("10:10").RowHeight = 24.75
$ws.Rows("11:11").RowHeight = 25.5
$ws.Rows("12:12").RowHeight = 16.5

# --- Update the totals row (now row 11) to reflect the new item's price. ---
$ws.Range("P11").Value = 52.88

# --- Update the generated timestamp shown in the footer. ---
$ws.Cells.Replace("Friday, 11 July, 2025 3:16 PM", "Friday, 11 July, 2025 4:00 PM")
